$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "name" column (C) entirely - display_id is now computed so a
# separate free-text name column is redundant; description shifts left.
$ws.Columns.Item(3).Delete()

# New explicit "key" values (cs-prefixed ids instead of bare slr/sll codes)
$ws.Range("A2").Value = "cs0002_slr0612"
$ws.Range("A3").Value = "cs0003_slr0613"
$ws.Range("A4").Value = "cs0004_sll0558"

# display_id is now derived from the key via formula instead of a static
# "{key}_codA_flat" string
$ws.Range("B2").Formula = '=CONCATENATE(A2,"_flat")'
$ws.Range("B3:B4").Formula = '=CONCATENATE(A3,"_flat")'

# description now references the {key} placeholder explicitly
$ws.Range("D2").Value = '"flattened version" of the design {key} suitable for visualization and genbank export'
$ws.Range("D3").Value = '"flattened version" of the design {key} suitable for visualization and genbank export'
$ws.Range("D4").Value = '"flattened version" of the design {key} suitable for visualization and genbank export'

# Match the column widths / selection left behind by the author's edit
$ws.Range("A1").EntireColumn.ColumnWidth = 13.28515625
$ws.Range("B1").EntireColumn.ColumnWidth = 19.85546875
$ws.Range("C1").EntireColumn.ColumnWidth = 25.7109375

$ws.Range("B2").Select()
